$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 85,2
$arr[0,0] = 'Cluster name'
$arr[0,1] = 'Active cases'
$arr[1,0] = '574 Plummer Street Building B1 & B2 Port Melbourne'
$arr[1,1] = 7
$arr[2,0] = '7 Chefs Fawkner'
$arr[2,1] = 5
$arr[3,0] = 'A1 Bakery Brunswick'
$arr[3,1] = 6
$arr[4,0] = 'Acquire BPO Southbank'
$arr[4,1] = 31
$arr[5,0] = 'Al Haj Halal Meats Glenroy'
$arr[5,1] = 72
$arr[6,0] = 'Al-Taqwa College Truganina'
$arr[6,1] = 23
$arr[7,0] = 'Amiga Montessori Craigieburn'
$arr[7,1] = 25
$arr[8,0] = 'Broadmeadows Medical Centre Broadmeadows'
$arr[8,1] = 5
$arr[9,0] = 'Budget Car and Truck Rentals Campbellfield'
$arr[9,1] = 5
$arr[10,0] = 'Campbellfield Heights Primary School Campbellfield'
$arr[10,1] = 6
$arr[11,0] = 'Can Panel Cambellfield'
$arr[11,1] = 6
$arr[12,0] = 'Cannie Road Construction Site Cannie'
$arr[12,1] = 7
$arr[13,0] = 'Caroline Springs Police Station'
$arr[13,1] = 7
$arr[14,0] = 'Cedars Medical Clinic Coburg'
$arr[14,1] = 43
$arr[15,0] = 'City of Hobsons Bay Community'
$arr[15,1] = 10
$arr[16,0] = 'City of Moreland Community'
$arr[16,1] = 9
$arr[17,0] = 'City of Wyndham Community'
$arr[17,1] = 5
$arr[18,0] = 'Classy Cabinets and Kitchens Craigieburn'
$arr[18,1] = 16
$arr[19,0] = 'Coles Barkly Square Brunswick August'
$arr[19,1] = 5
$arr[20,0] = 'Coles Broadmeadows Central Shopping Centre'
$arr[20,1] = 11
$arr[21,0] = 'Coles Campbellfield Plaza Campbellfield'
$arr[21,1] = 10
$arr[22,0] = 'Coles Coburg North Village'
$arr[22,1] = 10
$arr[23,0] = 'Coles Coburg North Village August'
$arr[23,1] = 8
$arr[24,0] = 'Coles Greenvale Shopping Centre'
$arr[24,1] = 6
$arr[25,0] = 'Coles Pakenham Place Shopping Centre'
$arr[25,1] = 6
$arr[26,0] = 'Coles Roxburgh Village Roxburgh Park'
$arr[26,1] = 17
$arr[27,0] = 'Community Kids Meadow Heights'
$arr[27,1] = 19
$arr[28,0] = 'Concept Caravans Campbellfield'
$arr[28,1] = 5
$arr[29,0] = 'Costco Wholesale Epping'
$arr[29,1] = 20
$arr[30,0] = 'Crossroads Logistics Sunshine North'
$arr[30,1] = 5
$arr[31,0] = 'Croydon Orthodontics'
$arr[31,1] = 7
$arr[32,0] = 'Elite Smart Community Care Campbellfield'
$arr[32,1] = 5
$arr[33,0] = 'Fitzroy Community School Fitzroy North'
$arr[33,1] = 52
$arr[34,0] = 'Glenroy West Primary School'
$arr[34,1] = 7
$arr[35,0] = 'Gloria Jeans Coffees Broadmeadows Central'
$arr[35,1] = 5
$arr[36,0] = 'Health Care Providers Association South Melbourne'
$arr[36,1] = 13
$arr[37,0] = 'IGA Meadow Heights Shopping Centre Meadow Heights'
$arr[37,1] = 6
$arr[38,0] = 'Ilim College Glenroy Campus Hadfield'
$arr[38,1] = 20
$arr[39,0] = 'Ilim College Kiewa Campus Boys Secondary Dallas'
$arr[39,1] = 6
$arr[40,0] = 'Ilim Learning Sanctuary Glenroy'
$arr[40,1] = 14
$arr[41,0] = 'Industrial Galvanizers Valmont Coatings Campbellfield'
$arr[41,1] = 13
$arr[42,0] = 'Islamic College of Melbourne Tarneit'
$arr[42,1] = 11
$arr[43,0] = 'KFC Fawkner'
$arr[43,1] = 20
$arr[44,0] = 'Kasr Sweets Coolaroo'
$arr[44,1] = 15
$arr[45,0] = 'Kids House Early Learning Cheltenham'
$arr[45,1] = 10
$arr[46,0] = 'Learning Nest Early Learning Centre Meadow Heights'
$arr[46,1] = 10
$arr[47,0] = 'Level Crossing Removal Project Lilydale Construction Site John Street'
$arr[47,1] = 8
$arr[48,0] = 'Louis Vuitton Melbourne Crown Southbank'
$arr[48,1] = 5
$arr[49,0] = 'Malvern Health and Fitness Clinic Malvern'
$arr[49,1] = 6
$arr[50,0] = 'Melbourne Metropolitan Remand Centre Ravenhall'
$arr[50,1] = 6
$arr[51,0] = 'Melbourne Truck Repairs Campbellfield'
$arr[51,1] = 6
$arr[52,0] = 'Melbourne West Police Station Docklands'
$arr[52,1] = 6
$arr[53,0] = 'Mercy Hospital for Women Heidelberg'
$arr[53,1] = 6
$arr[54,0] = 'Montessori Beginnings Greenvale'
$arr[54,1] = 6
$arr[55,0] = 'MyCentre Childcare Broadmeadows'
$arr[55,1] = 21
$arr[56,0] = 'Newbury Child and Community Centre Craigieburn'
$arr[56,1] = 10
$arr[57,0] = 'Newport Football Club Altona North'
$arr[57,1] = 7
$arr[58,0] = 'Newport Gardens Early Years Centre Newport'
$arr[58,1] = 5
$arr[59,0] = 'Nido Early School Moonee Ponds'
$arr[59,1] = 18
$arr[60,0] = 'Nino Early Learning Adventures Lalor'
$arr[60,1] = 5
$arr[61,0] = 'Nino Early Learning Adventures Lalor'
$arr[61,1] = 5
$arr[62,0] = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'
$arr[62,1] = 47
$arr[63,0] = 'OnQ Plumbing and Excavations Craigieburn'
$arr[63,1] = 8
$arr[64,0] = 'Oporto Coolaroo'
$arr[64,1] = 14
$arr[65,0] = 'Paisley Park Early Learning Centre Bundoora'
$arr[65,1] = 8
$arr[66,0] = 'Panorama Construction Site Whitehorse Rd Box Hill'
$arr[66,1] = 56
$arr[67,0] = 'Pearl Street Child Care Centre Glenroy'
$arr[67,1] = 6
$arr[68,0] = 'Salta Drive Construction Site Rangedale Drainage Altona North'
$arr[68,1] = 5
$arr[69,0] = 'Serco Mill Park'
$arr[69,1] = 8
$arr[70,0] = 'Southern Cross Station Crew Room Tissue Box Docklands'
$arr[70,1] = 5
$arr[71,0] = 'St Vincents Hospital Emergency Department Melbourne'
$arr[71,1] = 6
$arr[72,0] = 'Tek Foods Somerton'
$arr[72,1] = 9
$arr[73,0] = 'The Homestead Child and Family Centre Roxburgh Park'
$arr[73,1] = 18
$arr[74,0] = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'
$arr[74,1] = 5
$arr[75,0] = 'Tip Top Warehouse Dandenong'
$arr[75,1] = 10
$arr[76,0] = 'Total Window Concepts Hoppers Crossing'
$arr[76,1] = 6
$arr[77,0] = 'Tunstall Fresh Tunstall Square Shopping Centre Doncaster East'
$arr[77,1] = 5
$arr[78,0] = 'Victorian Civil and Administrative Tribunal Melbourne'
$arr[78,1] = 5
$arr[79,0] = 'Werribee Mercy Hospital Emergency Department'
$arr[79,1] = 6
$arr[80,0] = 'Western Health Footscray Hospital Emergency Department'
$arr[80,1] = 6
$arr[81,0] = 'Western Health Sunshine Hospital Emergency Department'
$arr[81,1] = 5
$arr[82,0] = 'Who is Bunker Spreckels Cafe Elwood'
$arr[82,1] = 5
$arr[83,0] = 'Woolworths Greenvale Lakes Roxburgh Park'
$arr[83,1] = 14
$arr[84,0] = 'Yara Childcare Centre Truganina'
$arr[84,1] = 27

$ws.Range("A1:B85").Value = $arr